# Commit: "added md_task to the plan table sheet"
#
# The "plan" worksheet gets a new row inserted at row 25 containing the
# md_task entity ("md_task" / "task_id task_version_id"), pushing the
# existing rows 25-38 down to 26-39. The active/selected sheet moves from
# "cdm" to "plan", and the remembered selections on the "MetaIdent" and
# "plan" sheets are updated to reflect where the editor was last working.

$wb = $excel.ActiveWorkbook

# --- "plan" sheet: insert the new md_task row -------------------------
$plan = $wb.Worksheets.Item("plan")
$plan.Rows.Item(25).Insert()
$plan.Cells.Item(25, 1).Value = "md_task"
$plan.Cells.Item(25, 2).Value = "task_id task_version_id"

# --- "MetaIdent" sheet: remembered selection moves to A34:B34 ---------
$meta = $wb.Worksheets.Item("MetaIdent")
$meta.Range("A34:B34").Select()

# --- "plan" sheet becomes the active tab, selection on B25 ------------
$plan.Activate()
$plan.Range("B25").Select()
